$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = '30.499.91'
$ws.Cells.Item(2, 5).Value = '  +0.42%  '
$ws.Cells.Item(3, 4).Value = '2.107.22'
$ws.Cells.Item(3, 5).Value = '  +4.79%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).Value = "'330.36"
$ws.Cells.Item(5, 5).Value = '  +1.66%  '
$ws.Cells.Item(6, 5).Value = '  +0.10%  '
$ws.Cells.Item(7, 4).Value = "'0.5256"
$ws.Cells.Item(7, 5).Value = '  +2.41%  '
$ws.Cells.Item(8, 4).Value = "'0.4390"
$ws.Cells.Item(8, 5).Value = '  +3.08%  '
$ws.Cells.Item(9, 4).Value = "'0.08885"
$ws.Cells.Item(9, 5).Value = '  +1.58%  '
$ws.Cells.Item(10, 4).Value = "'48.92"
$ws.Cells.Item(10, 5).Value = '  +12.26%  '
$ws.Cells.Item(11, 4).Value = "'1.167"
$ws.Cells.Item(11, 5).Value = '  +2.93%  '
$ws.Cells.Item(12, 4).Value = "'24.94"
$ws.Cells.Item(12, 5).Value = '  +2.12%  '
$ws.Cells.Item(13, 4).Value = '2.101.62'
$ws.Cells.Item(13, 5).Value = '  +4.21%  '
$ws.Cells.Item(14, 4).Value = "'6.755"
$ws.Cells.Item(14, 5).Value = '  +1.88%  '
$ws.Cells.Item(15, 4).Value = "'7.776"
$ws.Cells.Item(15, 5).Value = '  +4.33%  '
$ws.Cells.Item(16, 4).Value = "'96.64"
$ws.Cells.Item(16, 5).Value = '  +2.57%  '
$ws.Cells.Item(17, 4).Value = "'1.003"
$ws.Cells.Item(17, 5).Value = '  +0.08%  '
$ws.Cells.Item(18, 4).Value = "'0.00001131"
$ws.Cells.Item(18, 5).Value = '  +1.66%  '
$ws.Cells.Item(19, 4).Value = "'0.06643"
$ws.Cells.Item(19, 5).Value = '  +1.66%  '
$ws.Cells.Item(20, 4).Value = "'19.28"
$ws.Cells.Item(20, 5).Value = '  +2.35%  '
$ws.Cells.Item(21, 4).Value = "'1.002"
$ws.Cells.Item(21, 5).Value = '  +0.07%  '
$ws.Cells.Item(22, 5).Value = '  +1.82%  '
$ws.Cells.Item(23, 4).Value = '30.557.76'
$ws.Cells.Item(23, 5).Value = '  +0.37%  '
$ws.Cells.Item(24, 5).Value = '  +3.92%  '
$ws.Cells.Item(25, 4).Value = "'2.348"
$ws.Cells.Item(25, 5).Value = '  +4.23%  '
$ws.Cells.Item(26, 4).Value = '2.352.70'
$ws.Cells.Item(26, 5).Value = '  +4.44%  '
$ws.Cells.Item(27, 4).Value = "'22.49"
$ws.Cells.Item(27, 5).Value = '  +0.22%  '
$ws.Cells.Item(28, 4).Value = "'2.633"
$ws.Cells.Item(28, 5).Value = '  +7.30%  '
$ws.Cells.Item(29, 4).Value = "'162.06"
$ws.Cells.Item(29, 5).Value = '  +0.02%  '
$ws.Cells.Item(30, 4).Value = "'133.08"
$ws.Cells.Item(30, 5).Value = '  +1.31%  '
$ws.Cells.Item(31, 4).Value = "'1.228"
$ws.Cells.Item(31, 5).Value = '  +7.68%  '
$ws.Cells.Item(32, 5).Value = '  +1.86%  '
$ws.Cells.Item(33, 4).Value = "'1.693"
$ws.Cells.Item(33, 5).Value = '  +24.39%  '
$ws.Cells.Item(34, 4).Value = "'6.258"
$ws.Cells.Item(34, 5).Value = '  +2.84%  '
$ws.Cells.Item(35, 4).Value = "'3.892"
$ws.Cells.Item(35, 5).Value = '  +1.67%  '
$ws.Cells.Item(36, 5).Value = '  +11.14%  '
$ws.Cells.Item(37, 4).Value = "'0.02593"
$ws.Cells.Item(37, 5).Value = '  +2.22%  '
$ws.Cells.Item(38, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(38, 4).Value = "'5.528"
$ws.Cells.Item(38, 5).Value = '  +1.11%  '
$ws.Cells.Item(39, 2).Value = 'Hedera'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(39, 4).Value = "'0.06755"
$ws.Cells.Item(39, 5).Value = '  +1.18%  '
$ws.Cells.Item(40, 4).Value = "'12.79"
$ws.Cells.Item(40, 5).Value = '  +2.72%  '
$ws.Cells.Item(41, 4).Value = "'0.2289"
$ws.Cells.Item(41, 5).Value = '  +3.35%  '
$ws.Cells.Item(42, 4).Value = "'0.6941"
$ws.Cells.Item(42, 5).Value = '  +4.29%  '
$ws.Cells.Item(43, 4).Value = "'1.276"
$ws.Cells.Item(43, 5).Value = '  +2.90%  '
$ws.Cells.Item(44, 4).Value = "'1.001"
$ws.Cells.Item(44, 5).Value = '  +0.10%  '
$ws.Cells.Item(45, 4).Value = "'0.6438"
$ws.Cells.Item(45, 5).Value = '  +4.23%  '
$ws.Cells.Item(46, 4).Value = "'14.09"
$ws.Cells.Item(46, 5).Value = '  +3.07%  '
$ws.Cells.Item(47, 4).Value = "'2.231"
$ws.Cells.Item(47, 5).Value = '  +1.37%  '
$ws.Cells.Item(48, 4).Value = "'3.640"
$ws.Cells.Item(48, 5).Value = '  +0.18%  '
$ws.Cells.Item(49, 5).Value = '  -0.50%  '
$ws.Cells.Item(50, 5).Value = '  +10.04%  '
$ws.Cells.Item(51, 4).Value = "'82.91"
$ws.Cells.Item(51, 5).Value = '  +2.14%  '
